# Add a "Nama Siswa" (Student Name) column as the new first column of the
# sheet, shifting the existing "Alat Transportasi / Pekerjaan Orang Tua /
# Penghasilan Orang Tua / Jumlah Tanggungan / Pemilik KIP / Pemilik KPS"
# table one column to the right, then fill in the student names for the
# 10 existing data rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new column at A, pushing the current data (A:F) to (B:G).
$ws.Range("A1").EntireColumn.Insert()

# New header for the inserted column.
$ws.Range("A1").Value = "Nama Siswa"

# Student names for rows 2-11 (one per existing data row).
$namaSiswa = @("Ucup", "Ujank", "Umar", "Khabib", "Fadil", "Prakoso", "Fariz", "Trio", "Azel", "Dafa")

for ($i = 0; $i -lt $namaSiswa.Length; $i++) {
    $ws.Cells.Item($i + 2, 1).Value = $namaSiswa[$i]
}

# Match the saved selection/active cell from the edited workbook.
$null = $ws.Range("A11").Select()
